$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# Update the "Status" column (#512568 -> #512825) for the order rows
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 3).Value = "#512825"
}

# Remove the last order row (row 9): clear the product name entirely
# (content and formatting) and clear the quantity value while keeping
# its formatting.
$ws.Range("A9").Clear()
$ws.Range("B9").ClearContents()

# Update the active selection to D8
$ws.Range("D8").Select()
